$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Service section" (P/Q/R/S block) ---
# Row 4: task was "Merge code" / "All" -> rename task and reassign owner
$ws.Range("P4").Value = "Record Capstone"
$ws.Range("S4").Value = "Hoàng"

# Row 5: new task row added under the same block - copy formatting from row 4
# of the same block (same column styles) then fill in the values.
$ws.Range("P4:S4").Copy()
$ws.Range("P5:S5").PasteSpecial(-4122)
$ws.Range("P5").Value = "Merge code"
$ws.Range("Q5").Value = 45190
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = "All"

# --- "darkmode quote carousel section" (K/L/M/N block) ---
# M5 (% Build for the "Responsive" task) was left blank; mark as complete (100%),
# matching the filled-in style used by the sibling cells in that column (H5).
$ws.Range("H5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 1

# Row 5 grew a second wrapped line of content once P5:S5 were populated, so its
# height now matches the other two-line rows (1, 4, 6) in the sheet.
$ws.Rows.Item(5).RowHeight = 24.6

# Reflect the author's last clicked cell.
$ws.Range("L7").Select()

Write-Output "applied edits"
